$d = $word.ActiveDocument

# 1. Fix typo: "withing" -> "within"
$d.Content.Find.Execute(
    "specifically withing Waterville",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "specifically within Waterville", 2) | Out-Null

# 2. Insert a space: "empty/NA." -> "empty/ NA."
$d.Content.Find.Execute(
    "of cells that are empty/NA.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "of cells that are empty/ NA.", 2) | Out-Null

# 3. "experts familiar" -> "individuals  familiar" (note the double space,
#    matching the reflow of the run split in the source edit)
$d.Content.Find.Execute(
    "knowledge experts familiar with the BTM infestation issue.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "knowledge individuals  familiar with the BTM infestation issue.", 2) | Out-Null

# 4. "Filling in all NA values with median." -> "Filling in all NA values with the median."
$d.Content.Find.Execute(
    "Filling in all NA values with median.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Filling in all NA values with the median.", 2) | Out-Null

# 5. "Yes, The uses" -> "Yes, the uses" (capitalization fix)
$d.Content.Find.Execute(
    "Yes, The uses that I am aware of are for my class",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Yes, the uses that I am aware of are for my class", 2) | Out-Null
